$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ChromCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Rows 2-7: explicit unique values taken from the updated simulation run
$rowData = @{
    2  = @("111001011100001000100001001001", 0.8054944591353937,  0.0009282719617012632, 0.1861935271443584)
    3  = @("111001111100110110101111111110", 0.8198981795970176,  0.05242602698154674,   0.5640288929543333)
    4  = @("111001111100110110101111111110", 0.8198981795970176,  0.4296559397604132,    0.7404090244200711)
    5  = @("111001111100110110101111111110", 0.8198981795970176,  0.4296570337138129,    0.7808740650086972)
    6  = @("111001111100110110101111111110", 0.8198981795970176,  0.8198981795970176,    0.8198981795970177)
    7  = @("111001111100110110101111111110", 0.8198981795970176,  0.8198981795970176,    0.8198981795970177)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    Set-ChromCell "C$r" $vals[0]
    $ws.Range("D$r").Value = $vals[1]
    $ws.Range("E$r").Value = $vals[2]
    $ws.Range("F$r").Value = $vals[3]
}

# Rows 8-201: the run converged, every remaining row repeats the same stable values
$stableChrom = "111001111100110110101111111110"
$stableMax = 0.8198981795970176
$stableMin = 0.8198981795970176
$stableAvg = 0.8198981795970177

for ($r = 8; $r -le 201; $r++) {
    Set-ChromCell "C$r" $stableChrom
    $ws.Range("D$r").Value = $stableMax
    $ws.Range("E$r").Value = $stableMin
    $ws.Range("F$r").Value = $stableAvg
}
